# Remove display_name column from choices sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("choices")

# The "choices" sheet currently has columns:
#   A=list_name, B=name, C=display_name, D=label::language,
#   E=media::image::language, F=media::video::language, G=media::audio::language
# Delete column C (display_name), shifting everything left by one.
$ws.Columns.Item(3).Delete()

# Re-anchor the frozen panes: the freeze boundary was after column C (3 cols
# frozen); after deleting display_name it should be after column B (2 cols
# frozen) so the same logical columns (list_name, name) stay frozen.
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("C2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Select the (now-shifted) label::language column, matching the full-column
# selection left behind by the column deletion.
$ws.Columns.Item(3).Select()

# Restore the originally active sheet/tab.
$wb.Worksheets.Item("survey").Activate()
